$d = $word.ActiveDocument

# 1. Fix typo "foreigne" -> "foreign"
$d.Content.Find.Execute("foreigne", $false, $false, $false, $false, $false, $true, 1, $false, "foreign", 2)

# 2. Fix typo "wher the company" -> "where the company"
$d.Content.Find.Execute("wher the company", $false, $false, $false, $false, $false, $true, 1, $false, "where the company", 2)

# 3. Move the _GoBack bookmark so it sits right after "operating" in the
#    "Netflix recognizes ... where the company is operating" sentence,
#    instead of at the very end of the table cell.
$r = $d.Content
$r.Find.Execute("where the company is operating", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
